# Adds a new "bettermountHUD" mod entry at row 6 of the mod list, pushing
# the existing rows 6-28 down by one row (content + per-row formatting),
# which creates a new row 29 for the mod that used to be last (zoomify).
#
# Only columns A:G of the table are affected - column I (the "Legende"
# helper column used only in rows 3-6) is left completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 6..28 down to 7..29 (processing from the bottom up so a
# row's original contents are read before anything overwrites them).
# The destination is cleared first because PasteSpecial here leaves a
# cell's prior value alone when the source cell is blank.
for ($r = 29; $r -ge 7; $r--) {
    $src = $ws.Range("A" + ($r - 1) + ":G" + ($r - 1))
    $dst = $ws.Range("A" + $r + ":G" + $r)
    $dst.ClearContents()
    $src.Copy()
    $dst.PasteSpecial(-4122)   # xlPasteFormats
    $dst.PasteSpecial(-4163)   # xlPasteValues
}

$ws.Application.CutCopyMode = $false

# Row 6 becomes the new "bettermountHUD" entry. Its formatting (plain
# style, no "x" marker in column A) already matches what's needed, so
# only the mod-name value needs to change.
$ws.Range("A6").Value = ""
$ws.Range("B6").Value = "bettermountHUD"

# Restore the active-cell selection recorded in the saved workbook.
$ws.Range("I20").Select()
